$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append for distributor "ZV" (rows 62-77, columns B:G)
$newData = @(
    @("ZV", "2013 - 01", 4714, 191, 831072.88359999994, 11320.438),
    @("ZV", "2013 - 02", 3871, 129, 713605.09459999995, 7362.2052000000003),
    @("ZV", "2013 - 03", 3636, 129, 604788.12150000001, 6582.7501000000002),
    @("ZV", "2013 - 04", 4218, 173, 650873.3517, 10276.0908),
    @("ZV", "2013 - 05", 4168, 203, 612009.20880000002, 9432.9871999999996),
    @("ZV", "2013 - 06", 4580, 181, 732486.26870000002, 8934.3811000000005),
    @("ZV", "2013 - 07", 4817, 193, 770206.76729999995, 8869.8184000000001),
    @("ZV", "2013 - 08", 4818, 175, 809280.39919999999, 8307.0923000000003),
    @("ZV", "2013 - 09", 5187, 182, 879773.26650000003, 9466.6821),
    @("ZV", "2013 - 10", 4893, 177, 765758.50650000002, 8912.4038999999993),
    @("ZV", "2013 - 11", 4678, 189, 755234.00179999997, 10290.7323),
    @("ZV", "2013 - 12", 4132, 167, 698037.77060000005, 10579.5489),
    @("ZV", "2014 - 01", 4469, 170, 756851.92139999999, 11560.4812),
    @("ZV", "2014 - 02", 2510, 87, 400874.0601, 7578.5415000000003),
    @("ZV", "2014 - 03", 2301, 123, 411856.27360000001, 8555.5319999999992),
    @("ZV", "2014 - 04", 1490, 303, 311004.18589999998, 22993.235000000001)
)

$table = $ws.ListObjects.Item("Tabla1")

foreach ($rowValues in $newData) {
    $newRow = $table.ListRows.Add()
    $newRow.Range.Item(1, 1).Value = $rowValues[0]
    $newRow.Range.Item(1, 2).Value = $rowValues[1]
    $newRow.Range.Item(1, 3).Value = $rowValues[2]
    $newRow.Range.Item(1, 4).Value = $rowValues[3]
    $newRow.Range.Item(1, 5).Value = $rowValues[4]
    $newRow.Range.Item(1, 6).Value = $rowValues[5]
}
